$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate rows: 8 -> 10 -> 8 ---
$ws.Range("A10").Copy($ws.Range("A500"))
$ws.Range("A8").Copy($ws.Range("A10"))
$ws.Range("A500").Copy($ws.Range("A8"))
$ws.Range("A500").ClearContents()
$ws.Range("B10").Copy($ws.Range("B500"))
$ws.Range("B8").Copy($ws.Range("B10"))
$ws.Range("B500").Copy($ws.Range("B8"))
$ws.Range("B500").ClearContents()
$ws.Range("D10").Copy($ws.Range("D500"))
$ws.Range("D8").Copy($ws.Range("D10"))
$ws.Range("D500").Copy($ws.Range("D8"))
$ws.Range("D500").ClearContents()
$ws.Range("E10").Copy($ws.Range("E500"))
$ws.Range("E8").Copy($ws.Range("E10"))
$ws.Range("E500").Copy($ws.Range("E8"))
$ws.Range("E500").ClearContents()
$ws.Range("F10").Copy($ws.Range("F500"))
$ws.Range("F8").Copy($ws.Range("F10"))
$ws.Range("F500").Copy($ws.Range("F8"))
$ws.Range("F500").ClearContents()
$ws.Range("G10").Copy($ws.Range("G500"))
$ws.Range("G8").Copy($ws.Range("G10"))
$ws.Range("G500").Copy($ws.Range("G8"))
$ws.Range("G500").ClearContents()
$ws.Range("H10").Copy($ws.Range("H500"))
$ws.Range("H8").Copy($ws.Range("H10"))
$ws.Range("H500").Copy($ws.Range("H8"))
$ws.Range("H500").ClearContents()
$ws.Range("I10").Copy($ws.Range("I500"))
$ws.Range("I8").Copy($ws.Range("I10"))
$ws.Range("I500").Copy($ws.Range("I8"))
$ws.Range("I500").ClearContents()
$ws.Range("M8").Copy($ws.Range("M10"))
$ws.Range("M8").ClearContents()
$ws.Range("P10").Copy($ws.Range("P500"))
$ws.Range("P8").Copy($ws.Range("P10"))
$ws.Range("P500").Copy($ws.Range("P8"))
$ws.Range("P500").ClearContents()
$ws.Range("Q10").Copy($ws.Range("Q500"))
$ws.Range("Q8").Copy($ws.Range("Q10"))
$ws.Range("Q500").Copy($ws.Range("Q8"))
$ws.Range("Q500").ClearContents()
$ws.Range("R10").Copy($ws.Range("R500"))
$ws.Range("R8").Copy($ws.Range("R10"))
$ws.Range("R500").Copy($ws.Range("R8"))
$ws.Range("R500").ClearContents()
$ws.Range("S10").Copy($ws.Range("S500"))
$ws.Range("S8").Copy($ws.Range("S10"))
$ws.Range("S500").Copy($ws.Range("S8"))
$ws.Range("S500").ClearContents()
$ws.Range("T10").Copy($ws.Range("T500"))
$ws.Range("T8").Copy($ws.Range("T10"))
$ws.Range("T500").Copy($ws.Range("T8"))
$ws.Range("T500").ClearContents()
$ws.Range("U10").Copy($ws.Range("U500"))
$ws.Range("U8").Copy($ws.Range("U10"))
$ws.Range("U500").Copy($ws.Range("U8"))
$ws.Range("U500").ClearContents()
$ws.Range("V10").Copy($ws.Range("V500"))
$ws.Range("V8").Copy($ws.Range("V10"))
$ws.Range("V500").Copy($ws.Range("V8"))
$ws.Range("V500").ClearContents()
$ws.Range("W10").Copy($ws.Range("W500"))
$ws.Range("W8").Copy($ws.Range("W10"))
$ws.Range("W500").Copy($ws.Range("W8"))
$ws.Range("W500").ClearContents()
$ws.Range("Y10").Copy($ws.Range("Y500"))
$ws.Range("Y8").Copy($ws.Range("Y10"))
$ws.Range("Y500").Copy($ws.Range("Y8"))
$ws.Range("Y500").ClearContents()
$ws.Range("Z8").Copy($ws.Range("Z10"))
$ws.Range("Z8").ClearContents()
$ws.Range("AA10").Copy($ws.Range("AA500"))
$ws.Range("AA8").Copy($ws.Range("AA10"))
$ws.Range("AA500").Copy($ws.Range("AA8"))
$ws.Range("AA500").ClearContents()
$ws.Range("AB8").Copy($ws.Range("AB10"))
$ws.Range("AB8").ClearContents()
$ws.Range("AC8").Copy($ws.Range("AC10"))
$ws.Range("AC8").ClearContents()
$ws.Range("AD10").Copy($ws.Range("AD500"))
$ws.Range("AD8").Copy($ws.Range("AD10"))
$ws.Range("AD500").Copy($ws.Range("AD8"))
$ws.Range("AD500").ClearContents()
$ws.Range("AE10").Copy($ws.Range("AE500"))
$ws.Range("AE8").Copy($ws.Range("AE10"))
$ws.Range("AE500").Copy($ws.Range("AE8"))
$ws.Range("AE500").ClearContents()
$ws.Range("AG10").Copy($ws.Range("AG500"))
$ws.Range("AG8").Copy($ws.Range("AG10"))
$ws.Range("AG500").Copy($ws.Range("AG8"))
$ws.Range("AG500").ClearContents()
$ws.Range("AT10").Copy($ws.Range("AT500"))
$ws.Range("AT8").Copy($ws.Range("AT10"))
$ws.Range("AT500").Copy($ws.Range("AT8"))
$ws.Range("AT500").ClearContents()
$ws.Range("AW10").Copy($ws.Range("AW500"))
$ws.Range("AW8").Copy($ws.Range("AW10"))
$ws.Range("AW500").Copy($ws.Range("AW8"))
$ws.Range("AW500").ClearContents()
$ws.Range("AX10").Copy($ws.Range("AX500"))
$ws.Range("AX8").Copy($ws.Range("AX10"))
$ws.Range("AX500").Copy($ws.Range("AX8"))
$ws.Range("AX500").ClearContents()
$ws.Range("AY10").Copy($ws.Range("AY500"))
$ws.Range("AY8").Copy($ws.Range("AY10"))
$ws.Range("AY500").Copy($ws.Range("AY8"))
$ws.Range("AY500").ClearContents()

# --- Rotate rows: 11 -> 12 -> 11 ---
$ws.Range("A12").Copy($ws.Range("A500"))
$ws.Range("A11").Copy($ws.Range("A12"))
$ws.Range("A500").Copy($ws.Range("A11"))
$ws.Range("A500").ClearContents()
$ws.Range("B12").Copy($ws.Range("B500"))
$ws.Range("B11").Copy($ws.Range("B12"))
$ws.Range("B500").Copy($ws.Range("B11"))
$ws.Range("B500").ClearContents()
$ws.Range("D12").Copy($ws.Range("D500"))
$ws.Range("D11").Copy($ws.Range("D12"))
$ws.Range("D500").Copy($ws.Range("D11"))
$ws.Range("D500").ClearContents()
$ws.Range("E12").Copy($ws.Range("E500"))
$ws.Range("E11").Copy($ws.Range("E12"))
$ws.Range("E500").Copy($ws.Range("E11"))
$ws.Range("E500").ClearContents()
$ws.Range("F12").Copy($ws.Range("F500"))
$ws.Range("F11").Copy($ws.Range("F12"))
$ws.Range("F500").Copy($ws.Range("F11"))
$ws.Range("F500").ClearContents()
$ws.Range("G12").Copy($ws.Range("G500"))
$ws.Range("G11").Copy($ws.Range("G12"))
$ws.Range("G500").Copy($ws.Range("G11"))
$ws.Range("G500").ClearContents()
$ws.Range("H12").Copy($ws.Range("H500"))
$ws.Range("H11").Copy($ws.Range("H12"))
$ws.Range("H500").Copy($ws.Range("H11"))
$ws.Range("H500").ClearContents()
$ws.Range("I12").Copy($ws.Range("I500"))
$ws.Range("I11").Copy($ws.Range("I12"))
$ws.Range("I500").Copy($ws.Range("I11"))
$ws.Range("I500").ClearContents()
$ws.Range("M11").Copy($ws.Range("M12"))
$ws.Range("M11").ClearContents()
$ws.Range("P12").Copy($ws.Range("P500"))
$ws.Range("P11").Copy($ws.Range("P12"))
$ws.Range("P500").Copy($ws.Range("P11"))
$ws.Range("P500").ClearContents()
$ws.Range("Q12").Copy($ws.Range("Q500"))
$ws.Range("Q11").Copy($ws.Range("Q12"))
$ws.Range("Q500").Copy($ws.Range("Q11"))
$ws.Range("Q500").ClearContents()
$ws.Range("R12").Copy($ws.Range("R500"))
$ws.Range("R11").Copy($ws.Range("R12"))
$ws.Range("R500").Copy($ws.Range("R11"))
$ws.Range("R500").ClearContents()
$ws.Range("S12").Copy($ws.Range("S500"))
$ws.Range("S11").Copy($ws.Range("S12"))
$ws.Range("S500").Copy($ws.Range("S11"))
$ws.Range("S500").ClearContents()
$ws.Range("T12").Copy($ws.Range("T500"))
$ws.Range("T11").Copy($ws.Range("T12"))
$ws.Range("T500").Copy($ws.Range("T11"))
$ws.Range("T500").ClearContents()
$ws.Range("U12").Copy($ws.Range("U500"))
$ws.Range("U11").Copy($ws.Range("U12"))
$ws.Range("U500").Copy($ws.Range("U11"))
$ws.Range("U500").ClearContents()
$ws.Range("V12").Copy($ws.Range("V500"))
$ws.Range("V11").Copy($ws.Range("V12"))
$ws.Range("V500").Copy($ws.Range("V11"))
$ws.Range("V500").ClearContents()
$ws.Range("W12").Copy($ws.Range("W500"))
$ws.Range("W11").Copy($ws.Range("W12"))
$ws.Range("W500").Copy($ws.Range("W11"))
$ws.Range("W500").ClearContents()
$ws.Range("Y12").Copy($ws.Range("Y500"))
$ws.Range("Y11").Copy($ws.Range("Y12"))
$ws.Range("Y500").Copy($ws.Range("Y11"))
$ws.Range("Y500").ClearContents()
$ws.Range("Z11").Copy($ws.Range("Z12"))
$ws.Range("Z11").ClearContents()
$ws.Range("AA12").Copy($ws.Range("AA500"))
$ws.Range("AA11").Copy($ws.Range("AA12"))
$ws.Range("AA500").Copy($ws.Range("AA11"))
$ws.Range("AA500").ClearContents()
$ws.Range("AB11").Copy($ws.Range("AB12"))
$ws.Range("AB11").ClearContents()
$ws.Range("AC11").Copy($ws.Range("AC12"))
$ws.Range("AC11").ClearContents()
$ws.Range("AD12").Copy($ws.Range("AD500"))
$ws.Range("AD11").Copy($ws.Range("AD12"))
$ws.Range("AD500").Copy($ws.Range("AD11"))
$ws.Range("AD500").ClearContents()
$ws.Range("AE12").Copy($ws.Range("AE500"))
$ws.Range("AE11").Copy($ws.Range("AE12"))
$ws.Range("AE500").Copy($ws.Range("AE11"))
$ws.Range("AE500").ClearContents()
$ws.Range("AG12").Copy($ws.Range("AG500"))
$ws.Range("AG11").Copy($ws.Range("AG12"))
$ws.Range("AG500").Copy($ws.Range("AG11"))
$ws.Range("AG500").ClearContents()
$ws.Range("AT12").Copy($ws.Range("AT500"))
$ws.Range("AT11").Copy($ws.Range("AT12"))
$ws.Range("AT500").Copy($ws.Range("AT11"))
$ws.Range("AT500").ClearContents()
$ws.Range("AW12").Copy($ws.Range("AW500"))
$ws.Range("AW11").Copy($ws.Range("AW12"))
$ws.Range("AW500").Copy($ws.Range("AW11"))
$ws.Range("AW500").ClearContents()
$ws.Range("AX12").Copy($ws.Range("AX500"))
$ws.Range("AX11").Copy($ws.Range("AX12"))
$ws.Range("AX500").Copy($ws.Range("AX11"))
$ws.Range("AX500").ClearContents()
$ws.Range("AY12").Copy($ws.Range("AY500"))
$ws.Range("AY11").Copy($ws.Range("AY12"))
$ws.Range("AY500").Copy($ws.Range("AY11"))
$ws.Range("AY500").ClearContents()

# --- Rotate rows: 16 -> 17 -> 16 ---
$ws.Range("A17").Copy($ws.Range("A500"))
$ws.Range("A16").Copy($ws.Range("A17"))
$ws.Range("A500").Copy($ws.Range("A16"))
$ws.Range("A500").ClearContents()
$ws.Range("B17").Copy($ws.Range("B500"))
$ws.Range("B16").Copy($ws.Range("B17"))
$ws.Range("B500").Copy($ws.Range("B16"))
$ws.Range("B500").ClearContents()
$ws.Range("D17").Copy($ws.Range("D500"))
$ws.Range("D16").Copy($ws.Range("D17"))
$ws.Range("D500").Copy($ws.Range("D16"))
$ws.Range("D500").ClearContents()
$ws.Range("E17").Copy($ws.Range("E500"))
$ws.Range("E16").Copy($ws.Range("E17"))
$ws.Range("E500").Copy($ws.Range("E16"))
$ws.Range("E500").ClearContents()
$ws.Range("F17").Copy($ws.Range("F500"))
$ws.Range("F16").Copy($ws.Range("F17"))
$ws.Range("F500").Copy($ws.Range("F16"))
$ws.Range("F500").ClearContents()
$ws.Range("G17").Copy($ws.Range("G500"))
$ws.Range("G16").Copy($ws.Range("G17"))
$ws.Range("G500").Copy($ws.Range("G16"))
$ws.Range("G500").ClearContents()
$ws.Range("H17").Copy($ws.Range("H500"))
$ws.Range("H16").Copy($ws.Range("H17"))
$ws.Range("H500").Copy($ws.Range("H16"))
$ws.Range("H500").ClearContents()
$ws.Range("I17").Copy($ws.Range("I500"))
$ws.Range("I16").Copy($ws.Range("I17"))
$ws.Range("I500").Copy($ws.Range("I16"))
$ws.Range("I500").ClearContents()
$ws.Range("M16").Copy($ws.Range("M17"))
$ws.Range("M16").ClearContents()
$ws.Range("P17").Copy($ws.Range("P500"))
$ws.Range("P16").Copy($ws.Range("P17"))
$ws.Range("P500").Copy($ws.Range("P16"))
$ws.Range("P500").ClearContents()
$ws.Range("Q17").Copy($ws.Range("Q500"))
$ws.Range("Q16").Copy($ws.Range("Q17"))
$ws.Range("Q500").Copy($ws.Range("Q16"))
$ws.Range("Q500").ClearContents()
$ws.Range("R17").Copy($ws.Range("R500"))
$ws.Range("R16").Copy($ws.Range("R17"))
$ws.Range("R500").Copy($ws.Range("R16"))
$ws.Range("R500").ClearContents()
$ws.Range("S17").Copy($ws.Range("S500"))
$ws.Range("S16").Copy($ws.Range("S17"))
$ws.Range("S500").Copy($ws.Range("S16"))
$ws.Range("S500").ClearContents()
$ws.Range("T17").Copy($ws.Range("T500"))
$ws.Range("T16").Copy($ws.Range("T17"))
$ws.Range("T500").Copy($ws.Range("T16"))
$ws.Range("T500").ClearContents()
$ws.Range("U17").Copy($ws.Range("U500"))
$ws.Range("U16").Copy($ws.Range("U17"))
$ws.Range("U500").Copy($ws.Range("U16"))
$ws.Range("U500").ClearContents()
$ws.Range("V17").Copy($ws.Range("V500"))
$ws.Range("V16").Copy($ws.Range("V17"))
$ws.Range("V500").Copy($ws.Range("V16"))
$ws.Range("V500").ClearContents()
$ws.Range("W17").Copy($ws.Range("W500"))
$ws.Range("W16").Copy($ws.Range("W17"))
$ws.Range("W500").Copy($ws.Range("W16"))
$ws.Range("W500").ClearContents()
$ws.Range("Y17").Copy($ws.Range("Y500"))
$ws.Range("Y16").Copy($ws.Range("Y17"))
$ws.Range("Y500").Copy($ws.Range("Y16"))
$ws.Range("Y500").ClearContents()
$ws.Range("AA17").Copy($ws.Range("AA500"))
$ws.Range("AA16").Copy($ws.Range("AA17"))
$ws.Range("AA500").Copy($ws.Range("AA16"))
$ws.Range("AA500").ClearContents()
$ws.Range("AC16").Copy($ws.Range("AC17"))
$ws.Range("AC16").ClearContents()
$ws.Range("AD17").Copy($ws.Range("AD500"))
$ws.Range("AD16").Copy($ws.Range("AD17"))
$ws.Range("AD500").Copy($ws.Range("AD16"))
$ws.Range("AD500").ClearContents()
$ws.Range("AE17").Copy($ws.Range("AE500"))
$ws.Range("AE16").Copy($ws.Range("AE17"))
$ws.Range("AE500").Copy($ws.Range("AE16"))
$ws.Range("AE500").ClearContents()
$ws.Range("AG17").Copy($ws.Range("AG500"))
$ws.Range("AG16").Copy($ws.Range("AG17"))
$ws.Range("AG500").Copy($ws.Range("AG16"))
$ws.Range("AG500").ClearContents()
$ws.Range("AT17").Copy($ws.Range("AT500"))
$ws.Range("AT16").Copy($ws.Range("AT17"))
$ws.Range("AT500").Copy($ws.Range("AT16"))
$ws.Range("AT500").ClearContents()
$ws.Range("AW17").Copy($ws.Range("AW500"))
$ws.Range("AW16").Copy($ws.Range("AW17"))
$ws.Range("AW500").Copy($ws.Range("AW16"))
$ws.Range("AW500").ClearContents()
$ws.Range("AX17").Copy($ws.Range("AX500"))
$ws.Range("AX16").Copy($ws.Range("AX17"))
$ws.Range("AX500").Copy($ws.Range("AX16"))
$ws.Range("AX500").ClearContents()
$ws.Range("AY17").Copy($ws.Range("AY500"))
$ws.Range("AY16").Copy($ws.Range("AY17"))
$ws.Range("AY500").Copy($ws.Range("AY16"))
$ws.Range("AY500").ClearContents()

# --- Rotate rows: 20 -> 22 -> 20 ---
$ws.Range("A22").Copy($ws.Range("A500"))
$ws.Range("A20").Copy($ws.Range("A22"))
$ws.Range("A500").Copy($ws.Range("A20"))
$ws.Range("A500").ClearContents()
$ws.Range("B22").Copy($ws.Range("B500"))
$ws.Range("B20").Copy($ws.Range("B22"))
$ws.Range("B500").Copy($ws.Range("B20"))
$ws.Range("B500").ClearContents()
$ws.Range("D22").Copy($ws.Range("D500"))
$ws.Range("D20").Copy($ws.Range("D22"))
$ws.Range("D500").Copy($ws.Range("D20"))
$ws.Range("D500").ClearContents()
$ws.Range("E22").Copy($ws.Range("E500"))
$ws.Range("E20").Copy($ws.Range("E22"))
$ws.Range("E500").Copy($ws.Range("E20"))
$ws.Range("E500").ClearContents()
$ws.Range("F22").Copy($ws.Range("F500"))
$ws.Range("F20").Copy($ws.Range("F22"))
$ws.Range("F500").Copy($ws.Range("F20"))
$ws.Range("F500").ClearContents()
$ws.Range("G22").Copy($ws.Range("G500"))
$ws.Range("G20").Copy($ws.Range("G22"))
$ws.Range("G500").Copy($ws.Range("G20"))
$ws.Range("G500").ClearContents()
$ws.Range("H22").Copy($ws.Range("H500"))
$ws.Range("H20").Copy($ws.Range("H22"))
$ws.Range("H500").Copy($ws.Range("H20"))
$ws.Range("H500").ClearContents()
$ws.Range("I22").Copy($ws.Range("I500"))
$ws.Range("I20").Copy($ws.Range("I22"))
$ws.Range("I500").Copy($ws.Range("I20"))
$ws.Range("I500").ClearContents()
$ws.Range("M20").Copy($ws.Range("M22"))
$ws.Range("M20").ClearContents()
$ws.Range("P22").Copy($ws.Range("P500"))
$ws.Range("P20").Copy($ws.Range("P22"))
$ws.Range("P500").Copy($ws.Range("P20"))
$ws.Range("P500").ClearContents()
$ws.Range("Q22").Copy($ws.Range("Q500"))
$ws.Range("Q20").Copy($ws.Range("Q22"))
$ws.Range("Q500").Copy($ws.Range("Q20"))
$ws.Range("Q500").ClearContents()
$ws.Range("R22").Copy($ws.Range("R500"))
$ws.Range("R20").Copy($ws.Range("R22"))
$ws.Range("R500").Copy($ws.Range("R20"))
$ws.Range("R500").ClearContents()
$ws.Range("S22").Copy($ws.Range("S500"))
$ws.Range("S20").Copy($ws.Range("S22"))
$ws.Range("S500").Copy($ws.Range("S20"))
$ws.Range("S500").ClearContents()
$ws.Range("T22").Copy($ws.Range("T500"))
$ws.Range("T20").Copy($ws.Range("T22"))
$ws.Range("T500").Copy($ws.Range("T20"))
$ws.Range("T500").ClearContents()
$ws.Range("U22").Copy($ws.Range("U500"))
$ws.Range("U20").Copy($ws.Range("U22"))
$ws.Range("U500").Copy($ws.Range("U20"))
$ws.Range("U500").ClearContents()
$ws.Range("V22").Copy($ws.Range("V500"))
$ws.Range("V20").Copy($ws.Range("V22"))
$ws.Range("V500").Copy($ws.Range("V20"))
$ws.Range("V500").ClearContents()
$ws.Range("W22").Copy($ws.Range("W500"))
$ws.Range("W20").Copy($ws.Range("W22"))
$ws.Range("W500").Copy($ws.Range("W20"))
$ws.Range("W500").ClearContents()
$ws.Range("Y22").Copy($ws.Range("Y500"))
$ws.Range("Y20").Copy($ws.Range("Y22"))
$ws.Range("Y500").Copy($ws.Range("Y20"))
$ws.Range("Y500").ClearContents()
$ws.Range("Z22").Copy($ws.Range("Z500"))
$ws.Range("Z20").Copy($ws.Range("Z22"))
$ws.Range("Z500").Copy($ws.Range("Z20"))
$ws.Range("Z500").ClearContents()
$ws.Range("AA22").Copy($ws.Range("AA500"))
$ws.Range("AA20").Copy($ws.Range("AA22"))
$ws.Range("AA500").Copy($ws.Range("AA20"))
$ws.Range("AA500").ClearContents()
$ws.Range("AB22").Copy($ws.Range("AB500"))
$ws.Range("AB20").Copy($ws.Range("AB22"))
$ws.Range("AB500").Copy($ws.Range("AB20"))
$ws.Range("AB500").ClearContents()
$ws.Range("AC20").Copy($ws.Range("AC22"))
$ws.Range("AC20").ClearContents()
$ws.Range("AD22").Copy($ws.Range("AD500"))
$ws.Range("AD20").Copy($ws.Range("AD22"))
$ws.Range("AD500").Copy($ws.Range("AD20"))
$ws.Range("AD500").ClearContents()
$ws.Range("AE22").Copy($ws.Range("AE500"))
$ws.Range("AE20").Copy($ws.Range("AE22"))
$ws.Range("AE500").Copy($ws.Range("AE20"))
$ws.Range("AE500").ClearContents()
$ws.Range("AG22").Copy($ws.Range("AG500"))
$ws.Range("AG20").Copy($ws.Range("AG22"))
$ws.Range("AG500").Copy($ws.Range("AG20"))
$ws.Range("AG500").ClearContents()
$ws.Range("AT22").Copy($ws.Range("AT500"))
$ws.Range("AT20").Copy($ws.Range("AT22"))
$ws.Range("AT500").Copy($ws.Range("AT20"))
$ws.Range("AT500").ClearContents()
$ws.Range("AW22").Copy($ws.Range("AW500"))
$ws.Range("AW20").Copy($ws.Range("AW22"))
$ws.Range("AW500").Copy($ws.Range("AW20"))
$ws.Range("AW500").ClearContents()
$ws.Range("AX22").Copy($ws.Range("AX500"))
$ws.Range("AX20").Copy($ws.Range("AX22"))
$ws.Range("AX500").Copy($ws.Range("AX20"))
$ws.Range("AX500").ClearContents()
$ws.Range("AY22").Copy($ws.Range("AY500"))
$ws.Range("AY20").Copy($ws.Range("AY22"))
$ws.Range("AY500").Copy($ws.Range("AY20"))
$ws.Range("AY500").ClearContents()

# --- Rotate rows: 28 -> 29 -> 30 -> 28 ---
$ws.Range("A30").Copy($ws.Range("A500"))
$ws.Range("A29").Copy($ws.Range("A30"))
$ws.Range("A28").Copy($ws.Range("A29"))
$ws.Range("A500").Copy($ws.Range("A28"))
$ws.Range("A500").ClearContents()
$ws.Range("B30").Copy($ws.Range("B500"))
$ws.Range("B29").Copy($ws.Range("B30"))
$ws.Range("B28").Copy($ws.Range("B29"))
$ws.Range("B500").Copy($ws.Range("B28"))
$ws.Range("B500").ClearContents()
$ws.Range("D30").Copy($ws.Range("D500"))
$ws.Range("D29").Copy($ws.Range("D30"))
$ws.Range("D28").Copy($ws.Range("D29"))
$ws.Range("D500").Copy($ws.Range("D28"))
$ws.Range("D500").ClearContents()
$ws.Range("E30").Copy($ws.Range("E500"))
$ws.Range("E29").Copy($ws.Range("E30"))
$ws.Range("E28").Copy($ws.Range("E29"))
$ws.Range("E500").Copy($ws.Range("E28"))
$ws.Range("E500").ClearContents()
$ws.Range("F30").Copy($ws.Range("F500"))
$ws.Range("F29").Copy($ws.Range("F30"))
$ws.Range("F28").Copy($ws.Range("F29"))
$ws.Range("F500").Copy($ws.Range("F28"))
$ws.Range("F500").ClearContents()
$ws.Range("G30").Copy($ws.Range("G500"))
$ws.Range("G29").Copy($ws.Range("G30"))
$ws.Range("G28").Copy($ws.Range("G29"))
$ws.Range("G500").Copy($ws.Range("G28"))
$ws.Range("G500").ClearContents()
$ws.Range("H30").Copy($ws.Range("H500"))
$ws.Range("H29").Copy($ws.Range("H30"))
$ws.Range("H28").Copy($ws.Range("H29"))
$ws.Range("H500").Copy($ws.Range("H28"))
$ws.Range("H500").ClearContents()
$ws.Range("I30").Copy($ws.Range("I500"))
$ws.Range("I29").Copy($ws.Range("I30"))
$ws.Range("I28").Copy($ws.Range("I29"))
$ws.Range("I500").Copy($ws.Range("I28"))
$ws.Range("I500").ClearContents()
$ws.Range("M29").Copy($ws.Range("M30"))
$ws.Range("M29").ClearContents()
$ws.Range("M28").ClearContents()
$ws.Range("P30").Copy($ws.Range("P500"))
$ws.Range("P29").Copy($ws.Range("P30"))
$ws.Range("P28").Copy($ws.Range("P29"))
$ws.Range("P500").Copy($ws.Range("P28"))
$ws.Range("P500").ClearContents()
$ws.Range("Q30").Copy($ws.Range("Q500"))
$ws.Range("Q29").Copy($ws.Range("Q30"))
$ws.Range("Q28").Copy($ws.Range("Q29"))
$ws.Range("Q500").Copy($ws.Range("Q28"))
$ws.Range("Q500").ClearContents()
$ws.Range("R30").Copy($ws.Range("R500"))
$ws.Range("R29").Copy($ws.Range("R30"))
$ws.Range("R28").Copy($ws.Range("R29"))
$ws.Range("R500").Copy($ws.Range("R28"))
$ws.Range("R500").ClearContents()
$ws.Range("S30").Copy($ws.Range("S500"))
$ws.Range("S29").Copy($ws.Range("S30"))
$ws.Range("S28").Copy($ws.Range("S29"))
$ws.Range("S500").Copy($ws.Range("S28"))
$ws.Range("S500").ClearContents()
$ws.Range("T30").Copy($ws.Range("T500"))
$ws.Range("T29").Copy($ws.Range("T30"))
$ws.Range("T28").Copy($ws.Range("T29"))
$ws.Range("T500").Copy($ws.Range("T28"))
$ws.Range("T500").ClearContents()
$ws.Range("U30").Copy($ws.Range("U500"))
$ws.Range("U29").Copy($ws.Range("U30"))
$ws.Range("U28").Copy($ws.Range("U29"))
$ws.Range("U500").Copy($ws.Range("U28"))
$ws.Range("U500").ClearContents()
$ws.Range("V30").Copy($ws.Range("V500"))
$ws.Range("V29").Copy($ws.Range("V30"))
$ws.Range("V28").Copy($ws.Range("V29"))
$ws.Range("V500").Copy($ws.Range("V28"))
$ws.Range("V500").ClearContents()
$ws.Range("W30").Copy($ws.Range("W500"))
$ws.Range("W29").Copy($ws.Range("W30"))
$ws.Range("W28").Copy($ws.Range("W29"))
$ws.Range("W500").Copy($ws.Range("W28"))
$ws.Range("W500").ClearContents()
$ws.Range("Y30").Copy($ws.Range("Y500"))
$ws.Range("Y29").Copy($ws.Range("Y30"))
$ws.Range("Y28").Copy($ws.Range("Y29"))
$ws.Range("Y500").Copy($ws.Range("Y28"))
$ws.Range("Y500").ClearContents()
$ws.Range("Z30").Copy($ws.Range("Z500"))
$ws.Range("Z30").ClearContents()
$ws.Range("Z29").ClearContents()
$ws.Range("Z500").Copy($ws.Range("Z28"))
$ws.Range("Z500").ClearContents()
$ws.Range("AA30").Copy($ws.Range("AA500"))
$ws.Range("AA29").Copy($ws.Range("AA30"))
$ws.Range("AA28").Copy($ws.Range("AA29"))
$ws.Range("AA500").Copy($ws.Range("AA28"))
$ws.Range("AA500").ClearContents()
$ws.Range("AB30").Copy($ws.Range("AB500"))
$ws.Range("AB30").ClearContents()
$ws.Range("AB29").ClearContents()
$ws.Range("AB500").Copy($ws.Range("AB28"))
$ws.Range("AB500").ClearContents()
$ws.Range("AC29").Copy($ws.Range("AC30"))
$ws.Range("AC29").ClearContents()
$ws.Range("AC28").ClearContents()
$ws.Range("AD30").Copy($ws.Range("AD500"))
$ws.Range("AD29").Copy($ws.Range("AD30"))
$ws.Range("AD28").Copy($ws.Range("AD29"))
$ws.Range("AD500").Copy($ws.Range("AD28"))
$ws.Range("AD500").ClearContents()
$ws.Range("AE30").Copy($ws.Range("AE500"))
$ws.Range("AE29").Copy($ws.Range("AE30"))
$ws.Range("AE28").Copy($ws.Range("AE29"))
$ws.Range("AE500").Copy($ws.Range("AE28"))
$ws.Range("AE500").ClearContents()
$ws.Range("AG30").Copy($ws.Range("AG500"))
$ws.Range("AG29").Copy($ws.Range("AG30"))
$ws.Range("AG28").Copy($ws.Range("AG29"))
$ws.Range("AG500").Copy($ws.Range("AG28"))
$ws.Range("AG500").ClearContents()
$ws.Range("AT30").Copy($ws.Range("AT500"))
$ws.Range("AT29").Copy($ws.Range("AT30"))
$ws.Range("AT28").Copy($ws.Range("AT29"))
$ws.Range("AT500").Copy($ws.Range("AT28"))
$ws.Range("AT500").ClearContents()
$ws.Range("AW30").Copy($ws.Range("AW500"))
$ws.Range("AW29").Copy($ws.Range("AW30"))
$ws.Range("AW28").Copy($ws.Range("AW29"))
$ws.Range("AW500").Copy($ws.Range("AW28"))
$ws.Range("AW500").ClearContents()
$ws.Range("AX30").Copy($ws.Range("AX500"))
$ws.Range("AX29").Copy($ws.Range("AX30"))
$ws.Range("AX28").Copy($ws.Range("AX29"))
$ws.Range("AX500").Copy($ws.Range("AX28"))
$ws.Range("AX500").ClearContents()
$ws.Range("AY30").Copy($ws.Range("AY500"))
$ws.Range("AY29").Copy($ws.Range("AY30"))
$ws.Range("AY28").Copy($ws.Range("AY29"))
$ws.Range("AY500").Copy($ws.Range("AY28"))
$ws.Range("AY500").ClearContents()

# --- Rotate rows: 34 -> 33 -> 36 -> 34 ---
$ws.Range("A36").Copy($ws.Range("A500"))
$ws.Range("A33").Copy($ws.Range("A36"))
$ws.Range("A34").Copy($ws.Range("A33"))
$ws.Range("A500").Copy($ws.Range("A34"))
$ws.Range("A500").ClearContents()
$ws.Range("B36").Copy($ws.Range("B500"))
$ws.Range("B33").Copy($ws.Range("B36"))
$ws.Range("B34").Copy($ws.Range("B33"))
$ws.Range("B500").Copy($ws.Range("B34"))
$ws.Range("B500").ClearContents()
$ws.Range("D36").Copy($ws.Range("D500"))
$ws.Range("D33").Copy($ws.Range("D36"))
$ws.Range("D34").Copy($ws.Range("D33"))
$ws.Range("D500").Copy($ws.Range("D34"))
$ws.Range("D500").ClearContents()
$ws.Range("E36").Copy($ws.Range("E500"))
$ws.Range("E33").Copy($ws.Range("E36"))
$ws.Range("E34").Copy($ws.Range("E33"))
$ws.Range("E500").Copy($ws.Range("E34"))
$ws.Range("E500").ClearContents()
$ws.Range("F36").Copy($ws.Range("F500"))
$ws.Range("F33").Copy($ws.Range("F36"))
$ws.Range("F34").Copy($ws.Range("F33"))
$ws.Range("F500").Copy($ws.Range("F34"))
$ws.Range("F500").ClearContents()
$ws.Range("G36").Copy($ws.Range("G500"))
$ws.Range("G33").Copy($ws.Range("G36"))
$ws.Range("G34").Copy($ws.Range("G33"))
$ws.Range("G500").Copy($ws.Range("G34"))
$ws.Range("G500").ClearContents()
$ws.Range("H36").Copy($ws.Range("H500"))
$ws.Range("H33").Copy($ws.Range("H36"))
$ws.Range("H34").Copy($ws.Range("H33"))
$ws.Range("H500").Copy($ws.Range("H34"))
$ws.Range("H500").ClearContents()
$ws.Range("I36").Copy($ws.Range("I500"))
$ws.Range("I33").Copy($ws.Range("I36"))
$ws.Range("I34").Copy($ws.Range("I33"))
$ws.Range("I500").Copy($ws.Range("I34"))
$ws.Range("I500").ClearContents()
$ws.Range("M33").Copy($ws.Range("M36"))
$ws.Range("M33").ClearContents()
$ws.Range("M34").ClearContents()
$ws.Range("P36").Copy($ws.Range("P500"))
$ws.Range("P33").Copy($ws.Range("P36"))
$ws.Range("P34").Copy($ws.Range("P33"))
$ws.Range("P500").Copy($ws.Range("P34"))
$ws.Range("P500").ClearContents()
$ws.Range("Q36").Copy($ws.Range("Q500"))
$ws.Range("Q33").Copy($ws.Range("Q36"))
$ws.Range("Q34").Copy($ws.Range("Q33"))
$ws.Range("Q500").Copy($ws.Range("Q34"))
$ws.Range("Q500").ClearContents()
$ws.Range("R36").Copy($ws.Range("R500"))
$ws.Range("R33").Copy($ws.Range("R36"))
$ws.Range("R34").Copy($ws.Range("R33"))
$ws.Range("R500").Copy($ws.Range("R34"))
$ws.Range("R500").ClearContents()
$ws.Range("S36").Copy($ws.Range("S500"))
$ws.Range("S33").Copy($ws.Range("S36"))
$ws.Range("S34").Copy($ws.Range("S33"))
$ws.Range("S500").Copy($ws.Range("S34"))
$ws.Range("S500").ClearContents()
$ws.Range("T36").Copy($ws.Range("T500"))
$ws.Range("T33").Copy($ws.Range("T36"))
$ws.Range("T34").Copy($ws.Range("T33"))
$ws.Range("T500").Copy($ws.Range("T34"))
$ws.Range("T500").ClearContents()
$ws.Range("U36").Copy($ws.Range("U500"))
$ws.Range("U33").Copy($ws.Range("U36"))
$ws.Range("U34").Copy($ws.Range("U33"))
$ws.Range("U500").Copy($ws.Range("U34"))
$ws.Range("U500").ClearContents()
$ws.Range("V36").Copy($ws.Range("V500"))
$ws.Range("V33").Copy($ws.Range("V36"))
$ws.Range("V34").Copy($ws.Range("V33"))
$ws.Range("V500").Copy($ws.Range("V34"))
$ws.Range("V500").ClearContents()
$ws.Range("W36").Copy($ws.Range("W500"))
$ws.Range("W33").Copy($ws.Range("W36"))
$ws.Range("W34").Copy($ws.Range("W33"))
$ws.Range("W500").Copy($ws.Range("W34"))
$ws.Range("W500").ClearContents()
$ws.Range("Y36").Copy($ws.Range("Y500"))
$ws.Range("Y33").Copy($ws.Range("Y36"))
$ws.Range("Y34").Copy($ws.Range("Y33"))
$ws.Range("Y500").Copy($ws.Range("Y34"))
$ws.Range("Y500").ClearContents()
$ws.Range("Z36").Copy($ws.Range("Z500"))
$ws.Range("Z36").ClearContents()
$ws.Range("Z33").ClearContents()
$ws.Range("Z500").Copy($ws.Range("Z34"))
$ws.Range("Z500").ClearContents()
$ws.Range("AA36").Copy($ws.Range("AA500"))
$ws.Range("AA33").Copy($ws.Range("AA36"))
$ws.Range("AA34").Copy($ws.Range("AA33"))
$ws.Range("AA500").Copy($ws.Range("AA34"))
$ws.Range("AA500").ClearContents()
$ws.Range("AB36").Copy($ws.Range("AB500"))
$ws.Range("AB36").ClearContents()
$ws.Range("AB33").ClearContents()
$ws.Range("AB500").Copy($ws.Range("AB34"))
$ws.Range("AB500").ClearContents()
$ws.Range("AD36").Copy($ws.Range("AD500"))
$ws.Range("AD33").Copy($ws.Range("AD36"))
$ws.Range("AD34").Copy($ws.Range("AD33"))
$ws.Range("AD500").Copy($ws.Range("AD34"))
$ws.Range("AD500").ClearContents()
$ws.Range("AE36").Copy($ws.Range("AE500"))
$ws.Range("AE33").Copy($ws.Range("AE36"))
$ws.Range("AE34").Copy($ws.Range("AE33"))
$ws.Range("AE500").Copy($ws.Range("AE34"))
$ws.Range("AE500").ClearContents()
$ws.Range("AG36").Copy($ws.Range("AG500"))
$ws.Range("AG33").Copy($ws.Range("AG36"))
$ws.Range("AG34").Copy($ws.Range("AG33"))
$ws.Range("AG500").Copy($ws.Range("AG34"))
$ws.Range("AG500").ClearContents()
$ws.Range("AT36").Copy($ws.Range("AT500"))
$ws.Range("AT33").Copy($ws.Range("AT36"))
$ws.Range("AT34").Copy($ws.Range("AT33"))
$ws.Range("AT500").Copy($ws.Range("AT34"))
$ws.Range("AT500").ClearContents()
$ws.Range("AW36").Copy($ws.Range("AW500"))
$ws.Range("AW33").Copy($ws.Range("AW36"))
$ws.Range("AW34").Copy($ws.Range("AW33"))
$ws.Range("AW500").Copy($ws.Range("AW34"))
$ws.Range("AW500").ClearContents()
$ws.Range("AX36").Copy($ws.Range("AX500"))
$ws.Range("AX33").Copy($ws.Range("AX36"))
$ws.Range("AX34").Copy($ws.Range("AX33"))
$ws.Range("AX500").Copy($ws.Range("AX34"))
$ws.Range("AX500").ClearContents()
$ws.Range("AY36").Copy($ws.Range("AY500"))
$ws.Range("AY33").Copy($ws.Range("AY36"))
$ws.Range("AY34").Copy($ws.Range("AY33"))
$ws.Range("AY500").Copy($ws.Range("AY34"))
$ws.Range("AY500").ClearContents()
